$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.920.56'
$ws.Range("E2").Value = '  -2.15%  '
$ws.Range("D3").Value = '1.869.30'
$ws.Range("E3").Value = '  -2.24%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.ClearFormats()
$ws.Range("E4").Value = '  +0.13%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '311.97'
$c.ClearFormats()
$ws.Range("E5").Value = '  -1.04%  '
$ws.Range("E6").Value = '  +0.06%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4971'
$c.ClearFormats()
$ws.Range("E7").Value = '  -3.54%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3788'
$c.ClearFormats()
$ws.Range("E8").Value = '  -4.60%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.09014'
$c.ClearFormats()
$ws.Range("E9").Value = '  -7.46%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.114'
$c.ClearFormats()
$ws.Range("E10").Value = '  -3.23%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '41.66'
$c.ClearFormats()
$ws.Range("E11").Value = '  -1.54%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '6.301'
$c.ClearFormats()
$ws.Range("E12").Value = '  -3.63%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '20.64'
$c.ClearFormats()
$ws.Range("E13").Value = '  -2.65%  '
$ws.Range("D14").Value = '1.875.04'
$ws.Range("E14").Value = '  -1.98%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '7.212'
$c.ClearFormats()
$ws.Range("E15").Value = '  -3.83%  '
$ws.Range("E16").Value = '  +0.14%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.00001098'
$c.ClearFormats()
$ws.Range("E17").Value = '  -3.32%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '90.49'
$c.ClearFormats()
$ws.Range("E18").Value = '  -4.44%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06632'
$c.ClearFormats()
$ws.Range("E19").Value = '  -0.43%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '17.86'
$c.ClearFormats()
$ws.Range("E20").Value = '  -2.27%  '
$ws.Range("E21").Value = '  +0.15%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.103'
$c.ClearFormats()
$ws.Range("E22").Value = '  -3.55%  '
$ws.Range("D23").Value = '27.976.86'
$ws.Range("E23").Value = '  -2.15%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '11.34'
$c.ClearFormats()
$ws.Range("E24").Value = '  -1.49%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.289'
$c.ClearFormats()
$ws.Range("E25").Value = '  -1.06%  '
$ws.Range("D26").Value = '2.093.70'
$ws.Range("E26").Value = '  -1.71%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.509'
$c.ClearFormats()
$ws.Range("E27").Value = '  -6.38%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '157.99'
$c.ClearFormats()
$ws.Range("E28").Value = '  +0.09%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '20.66'
$c.ClearFormats()
$ws.Range("E29").Value = '  -2.87%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '126.00'
$c.ClearFormats()
$ws.Range("E30").Value = '  -2.29%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.1054'
$c.ClearFormats()
$ws.Range("E31").Value = '  -2.37%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.053'
$c.ClearFormats()
$ws.Range("E32").Value = '  -5.56%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.565'
$c.ClearFormats()
$ws.Range("E33").Value = '  -3.35%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '3.590'
$c.ClearFormats()
$ws.Range("E34").Value = '  -1.11%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '9.306'
$c.ClearFormats()
$ws.Range("E35").Value = '  -6.56%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.06524'
$c.ClearFormats()
$ws.Range("E36").Value = '  -3.97%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.02403'
$c.ClearFormats()
$ws.Range("E37").Value = '  -1.28%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.2181'
$c.ClearFormats()
$ws.Range("E38").Value = '  -1.68%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.268'
$c.ClearFormats()
$ws.Range("E39").Value = '  +6.47%  '
$ws.Range("E40").Value = '  -6.18%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '11.63'
$c.ClearFormats()
$ws.Range("E41").Value = '  -1.77%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.6354'
$c.ClearFormats()
$ws.Range("E42").Value = '  -2.00%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '4.886'
$c.ClearFormats()
$ws.Range("E43").Value = '  -4.14%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range("E44").Value = '  +0.05%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '13.15'
$c.ClearFormats()
$ws.Range("E45").Value = '  -3.00%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.5974'
$c.ClearFormats()
$ws.Range("E46").Value = '  -2.22%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.285'
$c.ClearFormats()
$ws.Range("E47").Value = '  -0.25%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '3.668'
$c.ClearFormats()
$ws.Range("E48").Value = '  -3.00%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.209'
$c.ClearFormats()
$ws.Range("E49").Value = '  +0.40%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.963'
$c.ClearFormats()
$ws.Range("E50").Value = '  -3.85%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '120.38'
$c.ClearFormats()
$ws.Range("E51").Value = '  -3.78%  '
